# ---------------------------------------------------------------------------
# Populate the "empty" 3-sheet workbook with real data, rename the sheets
# and restore the view/selection state exactly as it was left by the author.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# ---- rename the worksheets -------------------------------------------------
$ws1.Name = "Sheet1 - Text"
$ws2.Name = "Sheet2 - Numbers"
$ws3.Name = "Sheet3 - Formulas"

# ---- Sheet1 - Text ----------------------------------------------------------
$ws1.Range("A1").Value = "This is cell A1 in Sheet 1"
$ws1.Range("G5").Value = "This is cell G5"

# ---- Sheet2 - Numbers --------------------------------------------------------
# D1:D30 -> 1..30
# K1:K30 -> 1%..30%, formatted as a percentage (named "Percent" style)
# G5     -> shared text, same string as Sheet1!G5
$percents = @(0.01,0.02,0.03,0.04,0.05,0.06,0.07,0.08,0.09,0.10, `
              0.11,0.12,0.13,0.14,0.15,0.16,0.17,0.18,0.19,0.20, `
              0.21,0.22,0.23,0.24,0.25,0.26,0.27,0.28,0.29,0.30)

for ($i = 1; $i -le 30; $i++) {
    $ws2.Cells.Item($i, 4).Value = $i
    $ws2.Cells.Item($i, 11).Value = $percents[$i - 1]
}
$ws2.Range("K1:K30").Style = "Percent"
$ws2.Range("G5").Value = "This is cell G5"

# ---- Sheet3 - Formulas --------------------------------------------------------
$ws3.Range("D2").Formula = "='Sheet2 - Numbers'!D5"

# ---- selections / active cells, matching the saved view state --------------
[void]$ws1.Range("G6").Select()
[void]$ws2.Range("L2").Select()
[void]$ws3.Range("D3").Select()

# ---- active sheet (Sheet3 is the one left on-screen / tabSelected) ---------
[void]$ws3.Activate()
